$d = $word.ActiveDocument

function Replace-All($old, $new) {
    $rng = $d.Content
    $rng.Start = 0
    $count = 0
    while ($true) {
        $found = $rng.Find.Execute($old, $true, $false, $false, $false, $false,
                         $true, 1, $false, "", 0)
        if (-not $found) { break }
        $rng.Text = $new
        $count++
        if ($count -gt 50) { break }
    }
    if ($count -eq 0) {
        Write-Host "NOT FOUND: $old"
    }
}

Replace-All "Wamkelekile kwi ParentText" "Welcome to ParentText"

Replace-All "I-ParentText ifana nokuba nomhlobo okuxhasayo ecaleni kwakho, okubonisayo kulamahla ndenyuka okukhulisa umntwana wakho. " "ParentText is like having a supportive friend by your side, guiding you through the ups and downs of raising your child. "

Replace-All "Yenziwe ziingcali zase-Yunivesithi yaseKapa, i-Yunivesithi yaseOxford, i-Yunivesithi yase-Fort Hare, u-Clowns Without Borders South Africa kunye ne-Parenting for Lifelong Health, i-ParentText ivavanywe emhlabeni jikelele ukuqinisekisa ukuba inika olona ncedo lunokufumaneka. " "Created by a team of experts from the University of Cape Town, the University of Oxford, Fort Hare University, Clowns Without Borders South Africa, and Parenting for Lifelong Health, ParentText has been tested worldwide to ensure it offers the best help possible. "

Replace-All "Ndingu ______, umkhokheli wakho. Nangona ndibonakala ngathi ndingumntu, ndiyirobhothi eyenziwe yi-Parenting for Lifelong Health no-UNICEF ukuzokuxhasa kuhambo lwakho lokuba ngumzali. " "I am ______, your guide. Even though I might seem human, I am a robot created by Parenting for Lifelong Health and UNICEF to support you in your parenting journey. "

Replace-All "Masikhe sijonge ukuba isebenza kanjani i-ParentText. " "Let us see how ParentText works. "

Replace-All "I-ParentText ikunika izifundo ezintlanu zemihla ngemihla ukuphucula ubudlelwane bakho nomntwana wakho omncinci okanye omdala. Xa sele uzigqibile zontlanu izifundo zantsukuzonke, uzakufumana isatifiketi sokuba ngumzali omhle! " "ParentText offers 5 daily lessons to improve your relationship with your child or teen. Once you complete all 5 daily lessons, you will earn a positive parenting certificate! "

Replace-All "Phucula Ubudlelwane Bam noMntwana Wam" "Improve My Relationship with My Child or Teen"

Replace-All "Ukuchitha Ixesha Elikhethekileyo noMntwana Wam " "Spending One-on-one Time with My Child or Teen "

Replace-All "Ukumncoma " "Giving Praise "

Replace-All "Ukudala Inkqubo Yesiqhelo yeXesha Elikhethekileyo " "Creating a Routine for One-on-one Time "

Replace-All "Ukuqwalasela Iimvakalelo ngeliXesha likhethekileyo " "Noticing Feelings During One-on-one Time "

Replace-All " Ukuthoba Umoya Xa sinoXinzelelo " "Keeping Calm When We Are Stressed "

Replace-All "Ngoku, makhe sijonge ukuba sibonakala njani isifundo kwi ParentText. " "Now, let’s see what a lesson in ParentText looks like. "

Replace-All "Uzakufumana umyalezo yonke imihla okukhumbuzayo ukuba ugqibe isifundo sakho. Kwaye ukuba sikuphosile, kulungile! Usenokubuyela kwi ParentText nangaliphi na ixesha ukuze uqhubekele phambili nesifundo sakho." "You'll receive a daily notification to remind you to complete your lesson. And if you miss it, it is also okay! You can always return to ParentText anytime to catch up on your lesson."

Replace-All "Isifundo ngasinye siquka imibuzo, imifanekiso, Iingcebiso kunye nemidlalo emnandi ukuze uzame ukuyenza ekhaya nomntwana okanye nosapho lwakho." "Each lesson is a mix of quizzes, comics, tips, and a fun activity to try at home with your child or family."

Replace-All "Ukuba ukhe waxinga okanye ufuna uncedo, bhala MENYU okanye NCEDA ekupheleni kwesifundo sakho uzokufumana inkxaso eyongezelelweyo. " "If you are ever stuck or need help, type MENU or HELP at the end of your lessons to get more support. "

Replace-All "Xa ubhala NCEDA nanini na, ungafumana ulwazi ngezixhobo ezikhoyo ekuhlaleni ukumelana nobundlobongela bosapho, ubundlobongela ngokwesondo, impilo yengqondo, okanye nezinye iimeko zongxamiseko. " "When you type HELP anytime, you can get information about resources in your community to address family violence, sexual violence, mental health, or other emergencies. "

Replace-All "Iinkcukacha zakho zikhuselekile apha: Akukho nanye ekuzokwabelwana ngayo ngaphandle kwemvume yakho kwaye azizukuthengiswa ukwenza inzuzo. Le miyalezo oyithumelayo inoguqulelo oluntsokothileyo kwaye itshixelwe kwiseva ekhuselekileyo. " "Your information here is safe: Nothing will be shared without your permission and will not be sold for profit. The messages you send are encrypted and locked in a secure server. "

Replace-All "Khumbula, nabani na okwaziyo ukufikelela kwifowuni yakho xa ingatshixwanga angakwazi ukubona imiyalezo yakho. Ngoko ke, ukuba uthumela ulwazi olunobuzaza kwaye unexhala, cima imiyalezo kwifowuni yakho. " "Remember, anyone with access to your unlocked phone can view your messages. So, if you send sensitive information and are worried, delete the messages from your phone. "

Replace-All "UNCEDO" "HELP"

Replace-All "Ukuba khona kwakho apha kubonisa ukuba ukukhathalele kangakanani ukunika umntwana wakho eyona nkxaso. " "Being here shows how much you care about providing the best support for your child. "

Replace-All "Yinto oyenzayo nomntwana ezakwenza umahluko. " "It is what you do with your child that will really make a difference. "

Replace-All "I-ParentText iya kukubonelela ngeengcebiso ngezifundo ezizakunceda kubudlelwane bakho nomntwana wakho. Kuxhomekeke kuwe ukuba uyazisebenzisa ezingcebiso!" "ParentText will provide tips through lessons to help you with your relationship with your child. It is up to you to put these tips into practice!"

Replace-All "Enkosi kakhulu ngokumamela! Ungafikelela kulevidiyo nangaliphi na ixesha ngeMENU. Siyathemba uya kukonwabela ukusebenzisa i ParentText kwaye wenze lukhulu kuyo! " "Thank you so much for listening! You can access this video at any time via MENU. We hope you enjoy your ParentText journey and make the most out of it! "

Replace-All "Molo! Uziva njani njengangoku? Unayo imizuzwana engamashumi amathathu?" "Hi! How are you feeling right now? Do you have 30 seconds?"

Replace-All "Ngaphambi kokuba uqale kwi ParentText, masithathe ikhefu ekhawulezileyo kunye." "Before you get started in the ParentText programme, let's take a quick pause together."

Replace-All "Zama ukuThatha Ikhefu nanini na xa uziva unomsindo, unoxinzelelo, okanye ukhathazekile." "Try to Take a Pause whenever you feel angry, overwhelmed, stressed, or worried."

Replace-All "UngayiThatha Ikhefu nomntwana wakho omncinci okanye omdala!" "You can also Take a Pause with your child or teen!"

Replace-All "Thatha Ikhefu" "Take a Pause"

Replace-All "Hlala phantsi apho uzokhululeka khona ube sowuvala amehlo." "Sit down somewhere comfortable and close your eyes."

Replace-All "Phefumla nzuuulu." "Take a deeeeeeeep breath."

Replace-All "Wuve umoya ungena, uphuma, emzimbeni wakho." "Feel the air moving in, and out, of your body."

Replace-All "Wufake;" "In;"

Replace-All "wukhuphe;" "and out;"

Replace-All "Qwalasela ukuba uziva njani emzimbeni ngelixa uphefumlayo." "Notice how your body feels while you breathe."

Replace-All "Qwalasela ukuba uluva kweyiphi indawo uxinzelelo emzimbeni wakho." "Notice where you feel tension in your body."

Replace-All "Zama ukuyiphumza londawo." "Try to let it relax."

Replace-All "Xa sele ulungile, vula amehlo kwakhona." "When you are ready, open your eyes again."

Replace-All "Ngoku, qwalasela ukuba ingaba uziva ngokwahlukileyo kunangokuya" "Now, notice if you are feeling any differently than"

Replace-All "xa ubusaqala lomsebenzi." "when you started this activity."

Replace-All "Nokuphefumla nzulu kambalwa, okanye ukunxulumana nomhlaba, kungenza umehluko." "Even a few deep breaths, or connecting with the ground beneath you, can make a difference."

Replace-All "Thatha ikhefu nomntwana wakho omncinci okanye omdala!" "Take a pause with your child or teen!"
